$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($addr, $val)
    # Prefix with an apostrophe so Excel stores the literal text
    # (same as a user typing '0.4910 into a cell) instead of
    # re-parsing it into a Number, then clear the resulting
    # quote-prefix formatting so the cell keeps its original (default) style.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).ClearFormats()
}

# Row 2
Set-TextCell 'D2' '28.997.63'
Set-TextCell 'E2' '  -1.95%  '

# Row 3
Set-TextCell 'D3' '1.984.32'
Set-TextCell 'E3' '  -1.02%  '

# Row 4
Set-TextCell 'D4' '1.016'
Set-TextCell 'E4' '  +0.06%  '

# Row 5
Set-TextCell 'D5' '329.24'
Set-TextCell 'E5' '  -0.76%  '

# Row 6
Set-TextCell 'D6' '1.014'
Set-TextCell 'E6' '  +0.05%  '

# Row 7
Set-TextCell 'D7' '0.4910'
Set-TextCell 'E7' '  -2.37%  '

# Row 8
Set-TextCell 'D8' '0.4150'
Set-TextCell 'E8' '  -2.38%  '

# Row 9
Set-TextCell 'D9' '55.19'
Set-TextCell 'E9' '  +2.13%  '

# Row 10
Set-TextCell 'D10' '0.08826'
Set-TextCell 'E10' '  -3.64%  '

# Row 11
Set-TextCell 'D11' '1.081'
Set-TextCell 'E11' '  -3.97%  '

# Row 12
Set-TextCell 'D12' '2.075.05'
Set-TextCell 'E12' '  +2.35%  '

# Row 13
Set-TextCell 'D13' '22.67'
Set-TextCell 'E13' '  -3.87%  '

# Row 14
Set-TextCell 'D14' '7.868'
Set-TextCell 'E14' '  -3.00%  '

# Row 15
Set-TextCell 'E15' '  -2.97%  '

# Row 16
Set-TextCell 'E16' '  +0.20%  '

# Row 17
Set-TextCell 'D17' '91.38'
Set-TextCell 'E17' '  -4.68%  '

# Row 18
Set-TextCell 'D18' '0.00001096'
Set-TextCell 'E18' '  -2.50%  '

# Row 19
Set-TextCell 'D19' '0.06669'
Set-TextCell 'E19' '  -0.01%  '

# Row 20
Set-TextCell 'D20' '19.27'
Set-TextCell 'E20' '  -3.26%  '

# Row 21
Set-TextCell 'D21' '1.016'
Set-TextCell 'E21' '  +0.34%  '

# Row 22
Set-TextCell 'D22' '5.933'
Set-TextCell 'E22' '  -1.20%  '

# Row 23
Set-TextCell 'D23' '29.043.07'
Set-TextCell 'E23' '  -1.92%  '

# Row 24
Set-TextCell 'D24' '11.81'
Set-TextCell 'E24' '  -1.53%  '

# Row 25
Set-TextCell 'D25' '2.300'
Set-TextCell 'E25' '  +0.80%  '

# Row 26
Set-TextCell 'D26' '2.317.13'
Set-TextCell 'E26' '  +2.39%  '

# Row 27
Set-TextCell 'B27' 'EthereumClassic'
Set-TextCell 'C27' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 'D27' '20.64'
Set-TextCell 'E27' '  -0.80%  '

# Row 28
Set-TextCell 'B28' 'Monero'
Set-TextCell 'C28' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D28' '156.53'
Set-TextCell 'E28' '  -1.85%  '

# Row 29
Set-TextCell 'D29' '6.163'
Set-TextCell 'E29' '  -4.84%  '

# Row 30
Set-TextCell 'E30' '  -5.72%  '

# Row 31
Set-TextCell 'D31' '126.26'
Set-TextCell 'E31' '  -1.84%  '

# Row 32
Set-TextCell 'D32' '1.029'
Set-TextCell 'E32' '  -3.06%  '

# Row 33
Set-TextCell 'D33' '0.09806'
Set-TextCell 'E33' '  -1.53%  '

# Row 34
Set-TextCell 'D34' '1.507'
Set-TextCell 'E34' '  -5.13%  '

# Row 35
Set-TextCell 'D35' '5.792'
Set-TextCell 'E35' '  -1.68%  '

# Row 36
Set-TextCell 'D36' '3.745'
Set-TextCell 'E36' '  -1.67%  '

# Row 37
Set-TextCell 'D37' '0.02387'
Set-TextCell 'E37' '  -3.48%  '

# Row 38
Set-TextCell 'D38' '1.300'
Set-TextCell 'E38' '  -2.38%  '

# Row 39
Set-TextCell 'B39' 'FraxShare'
Set-TextCell 'C39' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D39' '8.950'
Set-TextCell 'E39' '  -6.83%  '

# Row 40
Set-TextCell 'B40' 'Hedera'
Set-TextCell 'C40' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D40' '0.06300'
Set-TextCell 'E40' '  -1.68%  '

# Row 41
Set-TextCell 'D41' '0.6428'
Set-TextCell 'E41' '  -2.46%  '

# Row 42
Set-TextCell 'D42' '11.41'
Set-TextCell 'E42' '  -3.44%  '

# Row 43
Set-TextCell 'D43' '1.015'

# Row 44
Set-TextCell 'D44' '0.1957'
Set-TextCell 'E44' '  -5.86%  '

# Row 45
Set-TextCell 'D45' '1.353'
Set-TextCell 'E45' '  +4.86%  '

# Row 46
Set-TextCell 'D46' '0.6134'
Set-TextCell 'E46' '  -3.84%  '

# Row 47
Set-TextCell 'D47' '13.17'
Set-TextCell 'E47' '  -3.86%  '

# Row 48
Set-TextCell 'D48' '2.135'
Set-TextCell 'E48' '  -3.57%  '

# Row 49
Set-TextCell 'B49' 'PancakeSwap'
Set-TextCell 'C49' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D49' '3.477'
Set-TextCell 'E49' '  -1.84%  '

# Row 50
Set-TextCell 'B50' 'BabyDogeCoin'
Set-TextCell 'C50' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 'D50' '0.00000000340'
Set-TextCell 'E50' '  +5.21%  '

# Row 51
Set-TextCell 'D51' '2.164'
Set-TextCell 'E51' '  +6.58%  '
